# Commit: "Coding more effects that can be applied to both enemies and
# player." While adding/reviewing entries in the shared Effects reference
# table further down the sheet, the author scrolled the window so row 54
# (the "average stats" / "Effects" section) became the first visible row,
# moving the view from topLeftCell A26 down to A54 - the active cell
# selection itself stays at A38, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate() | Out-Null

# Preserve the existing selection (A38) ...
$ws.Range("A38").Select() | Out-Null

# ... then scroll the window down so A54 becomes the new top-left visible cell.
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
